$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Current header order (row 1, columns C:F):
    #   C1 = normalize_group
    #   D1 = trajgroup_no_vary_q
    #   E1 = uniform_scaling_q
    #   F1 = variable_trajectory_group
    #
    # Target header order (row 1, columns C:F):
    #   C1 = variable_trajectory_group
    #   D1 = normalize_group
    #   E1 = trajgroup_no_vary_q
    #   F1 = uniform_scaling_q
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
